# Changed date and time fields for JGI app and verified database persistence
#
# The "survey" sheet documents the fields of the follow_map_time form. The
# FMT_FOL_date row (row 2) and the FMT_time row (row 4) had their "type"
# column corrected from the now-unsupported "date" / "time" field types to
# the generic "text" type.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Row 2 -> field FMT_FOL_date: type column C2 "date" -> "text"
$ws.Range("C2").Value = "text"

# Row 4 -> field FMT_time: type column C4 "time" -> "text"
$ws.Range("C4").Value = "text"

# The author's cursor ended up resting on C5 after the edit.
$ws.Range("C5").Select()

# Try to restore the recorded window geometry/zoom for parity with the
# author's session (best effort - some hosts do not persist window chrome).
$win = $excel.Windows.Item(1)
$win.Left = 3320
$win.Top = 180
$win.Width = 25600
$win.Height = 16060
